$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 1000
$ws.Range("B16").Value = 0.0571
$ws.Range("C16").Value = 0

$ws.Range("A17").Value = 2000
$ws.Range("B17").Value = 0.0521
$ws.Range("C17").Value = 0
